# plot mods fixes #89
# Re-order the rows of Table_1 so that, within each Year group, countries
# appear as Switzerland, Sweden, Spain. Concretely, the first and third
# data row of each 3-row year block are swapped while the middle
# ("Sweden") row is left untouched:
#   1890: rows 2 (Sweden) <-> 3 (Switzerland)
#   1918: rows 4 (Spain)  <-> 6 (Switzerland)     [row 5 Sweden unchanged]
#   1957: rows 7 (Spain)  <-> 9 (Switzerland)     [row 8 Sweden unchanged]
#   2020: rows 10 (Spain) <-> 12 (Switzerland)    [row 11 Sweden unchanged]

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch row far below the data (row 100) as temporary storage so
# that cell types (text vs number) and shared-string references are
# preserved exactly via Range.Copy instead of being re-parsed through
# .Value (which would mangle text like "12.0%" into a numeric percent).
function Swap-Rows($r1, $r2, $scratch) {
    $rng1 = $ws.Range("A" + $r1 + ":K" + $r1)
    $rng2 = $ws.Range("A" + $r2 + ":K" + $r2)
    $scratchRng = $ws.Range("A" + $scratch + ":K" + $scratch)

    $rng1.Copy($scratchRng)
    $rng2.Copy($rng1)
    $scratchRng.Copy($rng2)
    $scratchRng.Clear()
}

Swap-Rows 2 3 100
Swap-Rows 4 6 100
Swap-Rows 7 9 100
Swap-Rows 10 12 100
